$d = $word.ActiveDocument

# Locate the (index of the) paragraph that owns the "_GoBack" bookmark.
$bm = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm.Start
$bmIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $bmStart -and $p.Range.End -ge $bmStart) {
        $bmIndex = $i
    }
}

# Step 1: insert a whole new paragraph ("2023年3月3日") right before that paragraph.
$bmPara = $d.Paragraphs.Item($bmIndex)
[void]$bmPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($bmIndex)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>2023年3月3日</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$newPara.Range.InsertXML($xml1)

# Step 2: insert a run ("吃饭") right before the bookmark, inside its own (now shifted) paragraph.
$bm2 = $d.Bookmarks.Item("_GoBack")
$insPoint = $d.Range($bm2.Start, $bm2.Start)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>吃饭</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insPoint.InsertXML($xml2)

Write-Output "DONE"
